$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 2
Set-TextValue 2 4 '247.86'

# Row 4
Set-TextValue 4 4 '5.581'

# Row 5
Set-TextValue 5 4 '0.05623'

# Row 6
Set-TextValue 6 4 '3.403'

# Row 7
Set-TextValue 7 4 '6.482'

# Row 8
Set-TextValue 8 4 '0.8013'

# Row 9
Set-TextValue 9 4 '1.064'

# Row 10
$ws.Cells.Item(10, 2).Value = 'WazirX'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue 10 4 '0.1428'
$ws.Cells.Item(10, 5).Value = '9WazirXWRX'

# Row 11
$ws.Cells.Item(11, 2).Value = 'MandalaExchangeToken'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue 11 4 '0.07411'
$ws.Cells.Item(11, 5).Value = '10MandalaExchangeTokenMDX'

# Row 12
$ws.Cells.Item(12, 2).Value = 'LiechtensteinCryptoassetsExchange'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue 12 4 '0.03183'
$ws.Cells.Item(12, 5).Value = '11LiechtensteinCryptoassetsExchangeLCX'

# Row 13
$ws.Cells.Item(13, 2).Value = 'BitrueCoin'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue 13 4 '0.02974'
$ws.Cells.Item(13, 5).Value = '12BitrueCoinBTR'

# Row 14
$ws.Cells.Item(14, 2).Value = 'BitMartToken'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue 14 4 '0.09258'
$ws.Cells.Item(14, 5).Value = '13BitMartTokenBMX'

# Row 15
$ws.Cells.Item(15, 2).Value = 'BitForexToken'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue 15 4 '0.001659'
$ws.Cells.Item(15, 5).Value = '14BitForexTokenBF'

# Row 16
$ws.Cells.Item(16, 2).Value = 'CoinExToken'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
Set-TextValue 16 4 '0.04712'
$ws.Cells.Item(16, 5).Value = '15CoinExTokenCET'

# Row 17
$ws.Cells.Item(17, 2).Value = 'TigerCash'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue 17 4 '0.006262'
$ws.Cells.Item(17, 5).Value = '16TigerCashTCH'

# Row 18
$ws.Cells.Item(18, 2).Value = 'BitKan'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
Set-TextValue 18 4 '0.001055'
$ws.Cells.Item(18, 5).Value = '17BitKanKAN'

# Row 19
$ws.Cells.Item(19, 2).Value = 'HotbitToken'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
Set-TextValue 19 4 '0.003824'
$ws.Cells.Item(19, 5).Value = '18HotbitTokenHTB'

# Row 20
$ws.Cells.Item(20, 2).Value = 'NitroEx'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
Set-TextValue 20 4 '0.0001500'
$ws.Cells.Item(20, 5).Value = '19NitroExNTX'

# Row 21
$ws.Cells.Item(21, 2).Value = 'UpBots'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
Set-TextValue 21 4 '0.0004601'
$ws.Cells.Item(21, 5).Value = '20UpBotsUBXT'

# Row 22
$ws.Cells.Item(22, 2).Value = 'LEO'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue 22 4 '3.978'
$ws.Cells.Item(22, 5).Value = '21LEOLEO'

# Row 23
$ws.Cells.Item(23, 2).Value = 'BTSEToken'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue 23 4 '2.112'
$ws.Cells.Item(23, 5).Value = '22BTSETokenBTSE'

# Row 24
$ws.Cells.Item(24, 2).Value = 'One'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextValue 24 4 '0.01174'
$ws.Cells.Item(24, 5).Value = '23OneONEBestin24h'

# Row 26
$ws.Cells.Item(26, 5).Value = '25ProBitTokenPROB'

# Row 40
Set-TextValue 40 4 '0.04180'

# Row 41
Set-TextValue 41 4 '0.007182'

# Row 42
$ws.Cells.Item(42, 2).Value = 'CEJI'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
Set-TextValue 42 4 '0.003501'
$ws.Cells.Item(42, 5).Value = '41CEJICEJI'

# Row 43
$ws.Cells.Item(43, 2).Value = 'BKEXToken'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextValue 43 4 '0.1045'
$ws.Cells.Item(43, 5).Value = '42BKEXTokenBKK'

# Row 44
Set-TextValue 44 4 '0.008689'

# Row 45
Set-TextValue 45 4 '0.00005624'

# Row 47
Set-TextValue 47 4 '0.6802'

# Row 48
Set-TextValue 48 4 '0.02819'
